# Apply the cryptos.xlsx price/volume/ranking update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "67.467.20"
$ws.Range("E2").Value = "  -0.48%  "

$ws.Range("D3").Value = "2.506.40"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.93"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.69"
$ws.Range("E6").Value = "  -2.75%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D9").Value = "2.506.67"
$ws.Range("E9").Value = "  -0.79%  "

$ws.Range("E10").Value = "  -3.76%  "

$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  -2.01%  "

$ws.Range("B13").Value = "Toncoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.95"
$ws.Range("E13").Value = "  -4.17%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.03"
$ws.Range("E14").Value = "  -3.21%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.970.13"
$ws.Range("E15").Value = "  -0.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000174"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").Value = "67.338.40"
$ws.Range("E17").Value = "  -0.48%  "

$ws.Range("D18").Value = "2.489.43"
$ws.Range("E18").Value = "  -1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.70"
$ws.Range("E19").Value = "  +2.29%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.96"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.04"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("E22").Value = "  -3.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.47"
$ws.Range("E23").Value = "  -4.31%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.45"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("E26").Value = "  -6.32%  "

$ws.Range("E27").Value = "  -7.54%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.06"
$ws.Range("E28").Value = "  +5.91%  "

$ws.Range("D29").Value = "2.634.47"
$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("D30").Value = "0.0₃0936"
$ws.Range("E30").Value = "  -5.70%  "

$ws.Range("E31").Value = "  -1.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "519.67"
$ws.Range("E32").Value = "  -4.38%  "

$ws.Range("E33").Value = "  -1.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("E34").Value = "  -5.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  -1.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.17"
$ws.Range("E37").Value = "  -0.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.10"
$ws.Range("E38").Value = "  +1.33%  "

$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.60"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("E41").Value = "  -2.43%  "

$ws.Range("E42").Value = "  -3.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.333"
$ws.Range("E43").Value = "  -6.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("E45").Value = "  -0.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.25"
$ws.Range("E46").Value = "  -1.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.73"
$ws.Range("E47").Value = "  -1.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.64"
$ws.Range("E48").Value = "  -2.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.539"
$ws.Range("E49").Value = "  -4.14%  "

$ws.Range("D50").Value = "0.0₆0267"
$ws.Range("E50").Value = "  -3.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.67"
$ws.Range("E51").Value = "  -1.73%  "

